$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = 1.62
$ws.Range("P2").Value = 2.2
$ws.Range("Q2").Value = 3.1
$ws.Range("R2").Value = 1.36

# Row 5
$ws.Range("G5").Value = 3.2
$ws.Range("I5").Value = 2.4
$ws.Range("AR5").Value = 101

# Row 6
$ws.Range("Q6").Value = 2.06
$ws.Range("R6").Value = 1.84
$ws.Range("AG6").Value = 8.5
$ws.Range("AJ6").Value = 26
$ws.Range("AP6").Value = 23

# Row 7
$ws.Range("G7").Value = 4.1
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 1.9
$ws.Range("K7").Value = 2.05
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
$ws.Range("Q7").Value = 2.25
$ws.Range("R7").Value = 1.62
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.73
$ws.Range("W7").Value = 9.5
$ws.Range("AC7").Value = 8
$ws.Range("AJ7").Value = 15
$ws.Range("AK7").Value = 17
$ws.Range("AY7").Value = 23

# Row 10
$ws.Range("G10").Value = 2.2
$ws.Range("I10").Value = 3.5
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7
$ws.Range("X10").Value = 9.5
$ws.Range("AG10").Value = 9
$ws.Range("AH10").Value = 17
$ws.Range("AU10").Value = 8

# Row 13
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 10
$ws.Range("Q13").Value = 2.08
$ws.Range("R13").Value = 1.73
